$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the header formatting from G1 (bold/border/centered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the Save column values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
